# #327 Ajout des profils d'acces
# 1) Metadata sheet: bump the generation Date value.
$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# 2) Elements sheet: swap the two mapping columns (AK <-> AL), header included,
#    so "Mapping: Spécification métier vers l'extension ROR LocationEquipment"
#    moves before "Mapping: RIM Mapping".
$elements = $wb.Worksheets.Item("Elements")

$firstRow = 1
$lastRow = 17
$colAK = 37
$colAL = 38

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $akCell = $elements.Cells.Item($r, $colAK)
    $alCell = $elements.Cells.Item($r, $colAL)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    if ($akVal -ne $alVal) {
        $akCell.Value2 = $alVal
        $alCell.Value2 = $akVal
    }
}

# Column widths follow the (now swapped) content's best-fit widths.
$elements.Columns.Item($colAK).ColumnWidth = 74.17447916666667
$elements.Columns.Item($colAL).ColumnWidth = 24.147135416666668
